# Fix Training Data Issue
# The BF column ("Date") on this sheet held a malformed date string
# ("6-16-2012-13") for every data row. Correct it to the proper
# ISO-style date string "2013-06-16" for rows 2 through 31 (BF2:BF31),
# keeping the cell as plain text (not converting it into a real Excel
# date serial number/value).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 31; $row++) {
    $cell = $ws.Range("BF" + $row)
    # Prefix with an apostrophe so Excel stores this as literal text
    # instead of auto-converting the "YYYY-MM-DD" looking string into
    # a date value/serial number.
    $cell.Value = "'2013-06-16"
    # Remove the quote-prefix / text formatting that got applied so the
    # cell's style stays identical to its original (unstyled) state.
    $cell.ClearFormats()
}
